$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet1: "Overview"
# ============================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Drop all existing hyperlinks (will be re-created below with the
# refreshed display text, keeping the same underlying targets).
$ws1.Hyperlinks.Delete()

# Row 5 (".localization-config" row) collapses away - the handed-off
# file list now only has two real entries plus the config row.
$ws1.Rows.Item(5).Delete()

$ws1.Range("A2").Value = "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "25617f00-332c-4c35-a3c7-18b8487360b7.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/3e296d88-36a2-41da-bed1-0769611e3157.png", "", "", "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/68feb353-0183-491a-a576-84772a27d16c.png", "", "", "25617f00-332c-4c35-a3c7-18b8487360b7.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/7b5757aa-5efe-4c06-92ad-c607c6790e03.md", "", "", ".localization-config")

# ============================================================
# Sheet2: "zh-cn"
# ============================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(5).Delete()

$ws2.Range("A2").Value = "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.82d43905db1784f167d016f7b1e491e7708a0358.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-09 21:01:14"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "25617f00-332c-4c35-a3c7-18b8487360b7.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "25617f00-332c-4c35-a3c7-18b8487360b7.43438ea33b7dac54a24918443b8e4e5aa21ea5e6.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-09 21:01:14"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/3e296d88-36a2-41da-bed1-0769611e3157.png", "", "", "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83c4cd8ccd7dec16437e58bbb4346c9352e4cdcb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ac3eea108e1deeffec6050eac9913e6e2e39c7d7.png", "", "", "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.82d43905db1784f167d016f7b1e491e7708a0358.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/68feb353-0183-491a-a576-84772a27d16c.png", "", "", "25617f00-332c-4c35-a3c7-18b8487360b7.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83c4cd8ccd7dec16437e58bbb4346c9352e4cdcb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2169dcb74b5ab47500c108fb2e2cd54264a8959a.png", "", "", "25617f00-332c-4c35-a3c7-18b8487360b7.43438ea33b7dac54a24918443b8e4e5aa21ea5e6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/7b5757aa-5efe-4c06-92ad-c607c6790e03.md", "", "", ".localization-config")

# ============================================================
# Sheet3: "de-de"
# ============================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(5).Delete()

$ws3.Range("A2").Value = "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.82d43905db1784f167d016f7b1e491e7708a0358.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-09 21:01:24"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "25617f00-332c-4c35-a3c7-18b8487360b7.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "25617f00-332c-4c35-a3c7-18b8487360b7.43438ea33b7dac54a24918443b8e4e5aa21ea5e6.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-09 21:01:24"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/3e296d88-36a2-41da-bed1-0769611e3157.png", "", "", "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/135100f96d8f0a9bcef7af28e4080886e3d93609/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ac3eea108e1deeffec6050eac9913e6e2e39c7d7.png", "", "", "0ea6b364-551b-4f4a-86f7-b3f0a20abfbf.82d43905db1784f167d016f7b1e491e7708a0358.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/68feb353-0183-491a-a576-84772a27d16c.png", "", "", "25617f00-332c-4c35-a3c7-18b8487360b7.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/135100f96d8f0a9bcef7af28e4080886e3d93609/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2169dcb74b5ab47500c108fb2e2cd54264a8959a.png", "", "", "25617f00-332c-4c35-a3c7-18b8487360b7.43438ea33b7dac54a24918443b8e4e5aa21ea5e6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/625a01f0bf11bee0503e4dbd4ea40c006557de14/e2e/7b5757aa-5efe-4c06-92ad-c607c6790e03.md", "", "", ".localization-config")
